$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-S($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}
function Set-N($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = [double]$val
}

# Row 525
Set-S 525 1 "TRAINER_EDWIN_1"

# Row 526
Set-S 526 1 "species"
Set-S 526 2 "lvl"
Set-S 526 3 "iv"
Set-S 526 4 "heldItem"
Set-S 526 5 "moves"
Set-S 526 6 "ability"
Set-S 526 7 "shiny"

# Row 527
Set-S 527 1 "Turban"
Set-N 527 2 28
Set-S 527 5 "Razor Shell, Poison Fang, Iron Defense, Clamp"
Set-S 527 6 "Poison Touch"

# Row 529
Set-S 529 1 "TRAINER_JOSEPH"

# Row 530
Set-S 530 1 "species"
Set-S 530 2 "lvl"
Set-S 530 3 "iv"
Set-S 530 4 "heldItem"
Set-S 530 5 "moves"
Set-S 530 6 "ability"
Set-S 530 7 "shiny"

# Row 531
Set-S 531 1 "Galarian_Linoone"
Set-N 531 2 28

# Row 532
Set-S 532 1 "Electabuzz"
Set-N 532 2 26

# Row 534
Set-S 534 1 "TRAINER_ALYSSA"

# Row 535
Set-S 535 1 "species"
Set-S 535 2 "lvl"
Set-S 535 3 "iv"
Set-S 535 4 "heldItem"
Set-S 535 5 "moves"
Set-S 535 6 "ability"
Set-S 535 7 "shiny"

# Row 536
Set-S 536 1 "Doduo"
Set-N 536 2 27

# Row 537
Set-S 537 1 "Jolteon"
Set-N 537 2 29

# Row 539
Set-S 539 1 "TRAINER_EDWARD"

# Row 540
Set-S 540 1 "species"
Set-S 540 2 "lvl"
Set-S 540 3 "iv"
Set-S 540 4 "heldItem"
Set-S 540 5 "moves"
Set-S 540 6 "ability"
Set-S 540 7 "shiny"

# Row 541
Set-S 541 1 "Drowzee"
Set-N 541 2 26

# Row 542
Set-S 542 1 "Kadabra"
Set-N 542 2 28

# Row 544
Set-S 544 1 "TRAINER_DALE"

# Row 545
Set-S 545 1 "species"
Set-S 545 2 "lvl"
Set-S 545 3 "iv"
Set-S 545 4 "heldItem"
Set-S 545 5 "moves"
Set-S 545 6 "ability"
Set-S 545 7 "shiny"

# Row 546
Set-S 546 1 "Psyduck"
Set-N 546 2 27

# Row 547
Set-S 547 1 "Poliwhirl"
Set-N 547 2 28

# Row 549
Set-S 549 1 "TRAINER_JACLYN"

# Row 550
Set-S 550 1 "species"
Set-S 550 2 "lvl"
Set-S 550 3 "iv"
Set-S 550 4 "heldItem"
Set-S 550 5 "moves"
Set-S 550 6 "ability"
Set-S 550 7 "shiny"

# Row 551
Set-S 551 1 "Kirlia"
Set-N 551 2 28

# Row 552
Set-S 552 1 "Espeon"
Set-N 552 2 29

# Row 554
Set-S 554 1 "TRAINER_ABIGAIL_1"

# Row 555
Set-S 555 1 "species"
Set-S 555 2 "lvl"
Set-S 555 3 "iv"
Set-S 555 4 "heldItem"
Set-S 555 5 "moves"
Set-S 555 6 "ability"
Set-S 555 7 "shiny"

# Row 556
Set-S 556 1 "Ninjask"
Set-N 556 2 28

# Row 557
Set-S 557 1 "Fearow"
Set-N 557 2 29

# Row 559
Set-S 559 1 "TRAINER_ANTHONY"

# Row 560
Set-S 560 1 "species"
Set-S 560 2 "lvl"
Set-S 560 3 "iv"
Set-S 560 4 "heldItem"
Set-S 560 5 "moves"
Set-S 560 6 "ability"
Set-S 560 7 "shiny"

# Row 561
Set-S 561 1 "Electrode"
Set-N 561 2 29

# Row 562
Set-S 562 1 "Swellow"
Set-N 562 2 30

# Row 564
Set-S 564 1 "TRAINER_BENJAMIN_1"

# Row 565
Set-S 565 1 "species"
Set-S 565 2 "lvl"
Set-S 565 3 "iv"
Set-S 565 4 "heldItem"
Set-S 565 5 "moves"
Set-S 565 6 "ability"
Set-S 565 7 "shiny"

# Row 566
Set-S 566 1 "Growlithe"
Set-N 566 2 28

# Row 567
Set-S 567 1 "Linoone"
Set-N 567 2 30

# Row 569
Set-S 569 1 "TRAINER_JASMINE"

# Row 570
Set-S 570 1 "species"
Set-S 570 2 "lvl"
Set-S 570 3 "iv"
Set-S 570 4 "heldItem"
Set-S 570 5 "moves"
Set-S 570 6 "ability"
Set-S 570 7 "shiny"

# Row 571
Set-S 571 1 "Machoke"
Set-N 571 2 29

# Row 572
Set-S 572 1 "Metang"
Set-N 572 2 30

# Row 574
Set-S 574 1 "TRAINER_JACOB"

# Row 575
Set-S 575 1 "species"
Set-S 575 2 "lvl"
Set-S 575 3 "iv"
Set-S 575 4 "heldItem"
Set-S 575 5 "moves"
Set-S 575 6 "ability"
Set-S 575 7 "shiny"

# Row 576
Set-S 576 1 "Mankey"
Set-N 576 2 29

# Row 577
Set-S 577 1 "Metang"
Set-N 577 2 31

$null = $ws.Range("A579").Select()
